$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Angriffe")

# Fix A3 value (time) - small correction
$ws.Range("A3").Value = 44882.057733020833

# Add new rows 4 and 5
$ws.Range("A4:A5").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws.Range("A4:A5").HorizontalAlignment = -4131

$ws.Range("A4").Value = 44882.381963749998
$ws.Range("B4").Value = "Win"
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 515340
$ws.Range("E4").Value = 793494
$ws.Range("F4").Value = 3093

$ws.Range("A5").Value = 44882.472736634998
$ws.Range("B5").Value = "Win"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 193202
$ws.Range("E5").Value = 57778
$ws.Range("F5").Value = 294

# Select column A (whole column), matching the saved selection state
[void]$ws.Range("A1:A1048576").Select()
